# Update column F (dSF) values on Sheet1 to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -3
    3  = 3
    4  = -6
    7  = -7
    10 = 3
    11 = 2
    13 = 2
    15 = 4
    16 = 1
    18 = -5
    19 = -1
    20 = 1
    21 = 1
    22 = 1
    23 = 3
    24 = 2
    25 = 1
    27 = 4
    28 = -1
    29 = 6
    31 = -2
    32 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
